# "touppercase and xlxs fix for image links"
#
# Fixes two broken/incorrect image-filename references in the "Img url"
# column of the festival table:
#   - G6  (Stagecoach row): "../assets/images/stagecoachjpg"  -> "../assets/images/stagecoach.jpg"   (missing dot)
#   - G10 (Kaaboo row):     "../assets/images/kaboo.jpg"      -> "../assets/images/kaaboodelmar.jpg" (correct filename)
#
# Also restores the selection/active-cell state left in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Write G10 before G6 so the new shared-string entries land in the same
# order (kaaboodelmar.jpg, then stagecoach.jpg) as the authored workbook.
$ws.Range("G10").Value = "../assets/images/kaaboodelmar.jpg"
$ws.Range("G6").Value = "../assets/images/stagecoach.jpg"

# Update the view/selection to match the edited workbook.
$ws.Range("G6").Select()
